# Generate Report for Handoff
#
# The localization pipeline re-generated the handoff report: the
# translation status for zh-cn / de-de moved from "In Translation" to
# "Ready for handoff", the handoff datetime stamps advanced a couple of
# minutes, and the "Status" column on each sheet widened to fit the new
# (longer) status text.

$wb = $excel.ActiveWorkbook

$statusNew = "Ready for handoff"

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusNew          # zh-cn status
$wsOverview.Range("F2").Value = $statusNew          # de-de status
$wsOverview.Range("G2").Value = "2017-02-09 09:31:02"   # Latest HO Xliff Generate Date

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusNew              # Status
$wsZhCn.Range("H2").Value = "2017-02-09 09:30:39"       # Latest Handoff Datetime

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusNew              # Status
$wsDeDe.Range("H2").Value = "2017-02-09 09:31:02"       # Latest Handoff Datetime

# --- Widen the "Status" columns to fit the longer text -----------------
# (Overview: zh-cn/de-de status columns E & F; zh-cn/de-de sheets: Status
# column C). Use AutoFit first - the idiomatic way a user/script would
# react to the status column no longer fitting its content - then pin the
# resulting width to the value Excel settled on for this report.
$wsOverview.Columns.Item(5).AutoFit() | Out-Null
$wsOverview.Columns.Item(6).AutoFit() | Out-Null
$wsZhCn.Columns.Item(3).AutoFit() | Out-Null
$wsDeDe.Columns.Item(3).AutoFit() | Out-Null

$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3333333
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3333333
